$d = $word.ActiveDocument

$pairs = @(
    @("2025-07-01 Tuesday", "2025-07-02 Wednesday"),
    @("69×40=2760", "14×21=294"),
    @("77×32=2464", "58×88=5104"),
    @("86×24=2064", "21×64=1344"),
    @("87×12=1044", "53×84=4452"),
    @("53×81=4293", "17×48=816"),
    @("66×88=5808", "51×42=2142"),
    @("67×43=2881", "16×88=1408"),
    @("58×94=5452", "23×38=874"),
    @("59×71=4189", "38×93=3534"),
    @("31×68=2108", "65×39=2535"),
    @("16×36=576", "88×70=6160"),
    @("18×44=792", "32×18=576"),
    @("27×25=675", "14×25=350"),
    @("82×31=2542", "29×98=2842"),
    @("46×85=3910", "14×52=728"),
    @("75×17=1275", "91×31=2821"),
    @("66×71=4686", "17×15=255"),
    @("81×71=5751", "31×45=1395"),
    @("20×52=1040", "67×56=3752"),
    @("27×49=1323", "95×47=4465"),
    @("56×66=3696", "36×88=3168"),
    @("83×87=7221", "58×64=3712"),
    @("33×16=528", "33×12=396"),
    @("50×38=1900", "66×72=4752"),
    @("72×36=2592", "60×90=5400")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
